$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 507, pushing existing rows 507:595 down to 508:596
$ws.Rows.Item(507).Insert()

# Populate the newly inserted row 507 with the new data record
$ws.Range("A507").Value = 9
$ws.Range("B507").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C507").Value = "Metropolitana"
$ws.Range("D507").Value = 44798
$ws.Range("E507").Value = 13
$ws.Range("F507").Value = 100112024
$ws.Range("G507").Value = "Choclo"
$ws.Range("H507").Value = "Dulce o Americano"
$ws.Range("I507").Value = "Primera"
$ws.Range("J507").Value = 150
$ws.Range("K507").Value = 22000
$ws.Range("L507").Value = 25000
$ws.Range("M507").Value = 23600
$ws.Range("N507").Value = "$/malla 70 unidades"
$ws.Range("O507").Value = "Región de Arica y Parinacota"
$ws.Range("P507").Value = 337
$ws.Range("Q507").Value = 70
$ws.Range("R507").Value = "Hortaliza"
